# pareto-front + railroads + params
# Applies:
#  - configs sheet: new "Комментарий" column (C) + two new parameter rows
#    (pupilsPerCitizen, lackPenaltyCoefficient)
#  - limits sheet: new "Комментарий" column (D) with a comment per limit,
#    "park" limit flipped from included to excluded, and a new "railroad"
#    limit row
#  - cosmetic column width touch-ups on all three sheets

$wb = $excel.ActiveWorkbook

$wsProjects = $wb.Worksheets.Item("projects")
$wsConfigs  = $wb.Worksheets.Item("configs")
$wsLimits   = $wb.Worksheets.Item("limits")

# ---------------------------------------------------------------------
# configs sheet: add "Комментарий" header + two new rows
# ---------------------------------------------------------------------
$wsConfigs.Range("C1").Value = "Комментарий"

$wsConfigs.Range("A2").Value = "pupilsPerCitizen"
$wsConfigs.Range("B2").Value = 0.1
$wsConfigs.Range("C2").Value = "Доля учеников среди жителей"

$wsConfigs.Range("A3").Value = "lackPenaltyCoefficient"
$wsConfigs.Range("B3").Value = 5
$wsConfigs.Range("C3").Value = "Коэффициент для штрафа за нехватку учебных мест"

# ---------------------------------------------------------------------
# limits sheet: add "Комментарий" header + per-row comments
# ---------------------------------------------------------------------
$wsLimits.Range("D1").Value = "Комментарий"

$wsLimits.Range("D2").Value = "пляж"
$wsLimits.Range("D3").Value = "заправочные станции"
$wsLimits.Range("D4").Value = "промзоны"
$wsLimits.Range("D5").Value = "радиозоны"
$wsLimits.Range("D6").Value = "зоны плавки снега"
$wsLimits.Range("D7").Value = "транспортные узлы"
$wsLimits.Range("D8").Value = "вода"
$wsLimits.Range("D9").Value = "дороги"

# "park" row: now excluded (Включено = Нет) + comment
$wsLimits.Range("C10").Value = "Нет"
$wsLimits.Range("D10").Value = "парки"

$wsLimits.Range("D11").Value = "здания"

# new "railroad" limit row
$wsLimits.Range("A12").Value = "railroad"
$wsLimits.Range("B12").Value = 0.01
$wsLimits.Range("C12").Value = "Да"
$wsLimits.Range("D12").Value = "железные дороги"

# ---------------------------------------------------------------------
# Column width touch-ups (closest achievable values given the engine's
# pixel-quantized ColumnWidth -> stored-width conversion)
# ---------------------------------------------------------------------
$wsProjects.Columns.Item(2).ColumnWidth = 27.1669

$wsConfigs.Columns.Item(1).ColumnWidth = 24.6665
$wsConfigs.Columns.Item(2).ColumnWidth = 19.6664
$wsConfigs.Columns.Item(3).ColumnWidth = 55.1669

$wsLimits.Columns.Item(1).ColumnWidth = 18.0003
$wsLimits.Columns.Item(2).ColumnWidth = 27.8333
$wsLimits.Columns.Item(4).ColumnWidth = 29.4999

# ---------------------------------------------------------------------
# Selection / active sheet bookkeeping (matches the authored state)
# ---------------------------------------------------------------------
$wsConfigs.Range("C3").Select()
$wsLimits.Activate()
$wsLimits.Range("C12").Select()
